$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin (B) / Link (C) are plain text. Price (D) / Volume(1h) (E) look
# numeric (e.g. "1.000", "0.9998", "  -1.54%  ") so a leading apostrophe
# (the PowerShell escape `' inside a double-quoted string) forces Excel to
# store them as literal text -- matching the workbook's original inlineStr
# storage -- instead of silently coercing them to numbers.

$ws.Range("D2").Value = "`'30.228.89"
$ws.Range("E2").Value = "`'  -1.54%  "
$ws.Range("D3").Value = "`'1.849.16"
$ws.Range("E3").Value = "`'  -2.63%  "
$ws.Range("D4").Value = "`'0.9998"
$ws.Range("E4").Value = "`'  -0.09%  "
$ws.Range("D5").Value = "`'233.02"
$ws.Range("E5").Value = "`'  -2.59%  "
$ws.Range("D6").Value = "`'1.000"
$ws.Range("E6").Value = "`'  -0.01%  "
$ws.Range("E7").Value = "`'  -2.15%  "
$ws.Range("D8").Value = "`'0.2718"
$ws.Range("E8").Value = "`'  -4.59%  "
$ws.Range("D9").Value = "`'0.06374"
$ws.Range("E9").Value = "`'  -2.69%  "
$ws.Range("D10").Value = "`'1.835.39"
$ws.Range("E10").Value = "`'  -4.73%  "
$ws.Range("D11").Value = "`'0.07423"
$ws.Range("E11").Value = "`'  -0.44%  "
$ws.Range("D12").Value = "`'16.22"
$ws.Range("E12").Value = "`'  -2.97%  "
$ws.Range("D13").Value = "`'4.933"
$ws.Range("E13").Value = "`'  -3.40%  "
$ws.Range("D14").Value = "`'85.06"
$ws.Range("E14").Value = "`'  -3.49%  "
$ws.Range("D15").Value = "`'0.6273"
$ws.Range("E15").Value = "`'  -6.00%  "
$ws.Range("D16").Value = "`'30.171.30"
$ws.Range("E16").Value = "`'  -1.67%  "
$ws.Range("D17").Value = "`'0.9995"
$ws.Range("E17").Value = "`'  -0.05%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "`'12.60"
$ws.Range("E18").Value = "`'  -5.44%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "`'226.78"
$ws.Range("E19").Value = "`'  -2.06%  "
$ws.Range("D20").Value = "`'0.000007336"
$ws.Range("E20").Value = "`'  -3.64%  "
$ws.Range("D21").Value = "`'2.093.75"
$ws.Range("E21").Value = "`'  -4.99%  "
$ws.Range("D22").Value = "`'1.001"
$ws.Range("E22").Value = "`'  +0.00%  "
$ws.Range("D23").Value = "`'4.924"
$ws.Range("E23").Value = "`'  -6.89%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "`'5.924"
$ws.Range("E24").Value = "`'  -5.04%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "`'167.01"
$ws.Range("E25").Value = "`'  -1.79%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "`'9.216"
$ws.Range("E26").Value = "`'  -1.35%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "`'17.96"
$ws.Range("E27").Value = "`'  -4.30%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "`'1.873"
$ws.Range("E28").Value = "`'  -4.87%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "`'0.1020"
$ws.Range("E29").Value = "`'  +0.41%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "`'1.380"
$ws.Range("E30").Value = "`'  -1.68%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "`'4.111"
$ws.Range("E31").Value = "`'  -5.75%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "`'3.885"
$ws.Range("E32").Value = "`'  -3.63%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "`'0.04882"
$ws.Range("E33").Value = "`'  -4.74%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "`'1.155"
$ws.Range("E34").Value = "`'  -5.30%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "`'0.7101"
$ws.Range("E35").Value = "`'  -6.47%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "`'1.000"
$ws.Range("E36").Value = "`'  -0.18%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "`'2.699"
$ws.Range("E37").Value = "`'  -0.16%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "`'0.01841"
$ws.Range("E38").Value = "`'  -2.49%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "`'2.629"
$ws.Range("E39").Value = "`'  -1.20%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "`'0.9018"
$ws.Range("E40").Value = "`'  -2.21%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "`'1.950"
$ws.Range("E41").Value = "`'  -6.42%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "`'104.86"
$ws.Range("E42").Value = "`'  -1.98%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "`'0.9977"
$ws.Range("E43").Value = "`'  -0.74%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "`'0.4079"
$ws.Range("E44").Value = "`'  -5.36%  "
$ws.Range("D45").Value = "`'5.544"
$ws.Range("E45").Value = "`'  -3.38%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "`'7.065"
$ws.Range("E46").Value = "`'  -5.02%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "`'60.19"
$ws.Range("E47").Value = "`'  -6.75%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "`'0.1193"
$ws.Range("E48").Value = "`'  -6.55%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "`'8.601"
$ws.Range("E49").Value = "`'  -3.97%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "`'33.14"
$ws.Range("E50").Value = "`'  -2.28%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "`'1.384"
$ws.Range("E51").Value = "`'  -6.98%  "
